# ------------------------------------------------------------------
# "Chiffres COVID-19 Valais" daily update
# Updates quarantine contact counts (col O) / traveller counts (col P)
# for existing rows, appends 3 new days of data (rows 214-216 revised,
# row 217 filled in), and moves the frozen-pane / active selection to
# the new bottom of the data range.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column O/P revisions for rows 94-213 (retroactive count corrections) ---
$ws.Range("O94").Value = 32
$ws.Range("O95").Value = 33
$ws.Range("O96").Value = 36
$ws.Range("O97").Value = 40
$ws.Range("O98").Value = 46
$ws.Range("O99").Value = 48
$ws.Range("O100").Value = 53
$ws.Range("O101").Value = 70
$ws.Range("O102").Value = 66
$ws.Range("O103").Value = 82
$ws.Range("O104").Value = 95
$ws.Range("O105").Value = 98
$ws.Range("O106").Value = 108
$ws.Range("O107").Value = 108
$ws.Range("O108").Value = 124
$ws.Range("O109").Value = 115
$ws.Range("O110").Value = 123
$ws.Range("O111").Value = 116
$ws.Range("O112").Value = 124
$ws.Range("O113").Value = 114
$ws.Range("O114").Value = 117
$ws.Range("O115").Value = 153
$ws.Range("O116").Value = 168
$ws.Range("O117").Value = 180
$ws.Range("O118").Value = 162
$ws.Range("O119").Value = 162
$ws.Range("O120").Value = 157
$ws.Range("O121").Value = 165
$ws.Range("O122").Value = 172
$ws.Range("O123").Value = 165
$ws.Range("O124").Value = 159
$ws.Range("O125").Value = 127
$ws.Range("O126").Value = 124
$ws.Range("O127").Value = 121
$ws.Range("O128").Value = 129
$ws.Range("O129").Value = 128
$ws.Range("O130").Value = 124
$ws.Range("O131").Value = 118
$ws.Range("O132").Value = 109
$ws.Range("O133").Value = 92
$ws.Range("O134").Value = 83
$ws.Range("O135").Value = 74
$ws.Range("O136").Value = 70
$ws.Range("O137").Value = 52
$ws.Range("O138").Value = 74
$ws.Range("O139").Value = 79
$ws.Range("O140").Value = 82
$ws.Range("O141").Value = 84
$ws.Range("O142").Value = 81
$ws.Range("O143").Value = 82
$ws.Range("O144").Value = 81
$ws.Range("O145").Value = 79
$ws.Range("O146").Value = 62
$ws.Range("O147").Value = 65
$ws.Range("O148").Value = 47
$ws.Range("O149").Value = 49
$ws.Range("O150").Value = 56
$ws.Range("O151").Value = 55
$ws.Range("O152").Value = 52
$ws.Range("O153").Value = 55
$ws.Range("O154").Value = 59
$ws.Range("O155").Value = 87
$ws.Range("O156").Value = 139
$ws.Range("O157").Value = 156
$ws.Range("O158").Value = 150
$ws.Range("O159").Value = 164
$ws.Range("O160").Value = 169
$ws.Range("O161").Value = 159
$ws.Range("O162").Value = 156
$ws.Range("O163").Value = 169
$ws.Range("O164").Value = 184
$ws.Range("O165").Value = 178
$ws.Range("O166").Value = 166
$ws.Range("O167").Value = 193
$ws.Range("O168").Value = 192
$ws.Range("O169").Value = 208
$ws.Range("O170").Value = 215
$ws.Range("O171").Value = 230
$ws.Range("O172").Value = 220
$ws.Range("O173").Value = 195
$ws.Range("O174").Value = 241
$ws.Range("O175").Value = 258
$ws.Range("O176").Value = 274
$ws.Range("O177").Value = 284
$ws.Range("O178").Value = 278
$ws.Range("O179").Value = 285
$ws.Range("O180").Value = 360
$ws.Range("O181").Value = 365
$ws.Range("O182").Value = 372
$ws.Range("O183").Value = 413
$ws.Range("O184").Value = 433
$ws.Range("O185").Value = 459
$ws.Range("O186").Value = 445
$ws.Range("O187").Value = 426
$ws.Range("O188").Value = 357
$ws.Range("O189").Value = 316
$ws.Range("O190").Value = 296
$ws.Range("O191").Value = 260
$ws.Range("O192").Value = 222
$ws.Range("O193").Value = 232
$ws.Range("O194").Value = 241
$ws.Range("O195").Value = 212
$ws.Range("O196").Value = 213
$ws.Range("O197").Value = 224
$ws.Range("O198").Value = 213
$ws.Range("O199").Value = 227
$ws.Range("O200").Value = 233
$ws.Range("O201").Value = 256
$ws.Range("O202").Value = 294
$ws.Range("O203").Value = 288
$ws.Range("O204").Value = 330
$ws.Range("O205").Value = 384
$ws.Range("O206").Value = 407
$ws.Range("O207").Value = 426
$ws.Range("O208").Value = 419
$ws.Range("O209").Value = 384
$ws.Range("O210").Value = 398
$ws.Range("O211").Value = 403
$ws.Range("P211").Value = 440
$ws.Range("O212").Value = 426
$ws.Range("P212").Value = 476
$ws.Range("O213").Value = 446
$ws.Range("P213").Value = 481

# --- Rows 214-216: revised input figures ---
$ws.Range("C214").Value = 6
$ws.Range("N214").Value = 95
$ws.Range("O214").Value = 392
$ws.Range("P214").Value = 511

$ws.Range("C215").Value = 4
$ws.Range("N215").Value = 94
$ws.Range("O215").Value = 348
$ws.Range("P215").Value = 511

$ws.Range("C216").Value = 12
$ws.Range("I216").Value = 1
$ws.Range("N216").Value = 82
$ws.Range("O216").Value = 301
$ws.Range("P216").Value = 474

# --- Row 217: newly filled-in day (was blank) ---
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 1
$ws.Range("E217").Value = 1
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 8
$ws.Range("I217").Value = 0
$ws.Range("L217").Value = "0"
$ws.Range("M217").Value = "0"
$ws.Range("N217").Value = 75
$ws.Range("O217").Value = 262
$ws.Range("P217").Value = 409

# --- Recalculate so cached formula results (cumuls, totals) are fresh ---
$excel.CalculateFull()

# --- View state: frozen-pane anchor and active selection follow the new data ---
$ws.Range("H207").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("V221").Select()
